$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values to match repulled data
$ws.Range("F3").Value = -2
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 2
$ws.Range("F6").Value = 3
$ws.Range("F8").Value = 2
$ws.Range("F9").Value = -1
$ws.Range("F10").Value = -2
$ws.Range("F11").Value = 2
$ws.Range("F12").Value = -4
$ws.Range("F14").Value = 4
$ws.Range("F15").Value = 1
$ws.Range("F16").Value = -3
$ws.Range("F18").Value = -2
$ws.Range("F21").Value = 5
$ws.Range("F22").Value = -3
$ws.Range("F23").Value = 0
